# Commit: "Inclusão da União PEAD DN 20 como contratada."
#
# Adds a new "TROCA DE CAIXA DE PARADA" (322000) service-type entry to the
# "unitario" sheet's lookup table (inserted right before the existing
# "REPARO DE REDE DE AGUA NAO VISIVEL" row, pushing every row below it down
# by one), and brings the "unitario" tab back into focus (it had been left
# on "reposicao").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unitario")

# --- Insert the new catalog row at row 30 --------------------------------
# (everything that used to live at row 30 onward shifts down by one row)
$ws.Rows.Item(30).Insert()

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(30).RowHeight = 18.75

$ws.Range("A30").Value = "322000"
$ws.Range("B30").Value = "TROCA DE CAIXA DE PARADA"
$ws.Range("C30").Value = "Poco"

# --- Keep the (stale) _FilterDatabase defined name in sync ---------------
# It referenced unitario!$A$4:$C$53 and must now cover the extra row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "unitario!_FilterDatabase") {
        $n.RefersTo = "=unitario!`$A`$4:`$C`$54"
    }
}

# --- Switch the active tab back to "unitario" -----------------------------
# The workbook was last saved with "reposicao" focused; the author
# reviewed/edited "unitario" last, so it becomes the selected tab.
$ws.Activate()
